$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 0
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 2
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 2
$ws.Range("D12").Value = 4
$ws.Range("E12").Value = 3
$ws.Range("D13").Value = 5
$ws.Range("E13").Value = 4
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 0
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 2
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = 2
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = 5
$ws.Range("F17").Value = 3
$ws.Range("D18").Value = 11
$ws.Range("E18").Value = 6
$ws.Range("F18").Value = 4
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = 9
$ws.Range("F19").Value = 4
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 15
$ws.Range("E20").Value = 10
$ws.Range("F20").Value = 2
$ws.Range("C21").Value = 3
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 11
$ws.Range("F21").Value = 3
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 21
$ws.Range("E22").Value = 13
$ws.Range("F22").Value = 4
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 24
$ws.Range("E23").Value = 14
$ws.Range("F23").Value = 5
$ws.Range("B24").Value = 0
$ws.Range("C24").Value = 5
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = 14
$ws.Range("F24").Value = 6
$ws.Range("B25").Value = 0
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 27
$ws.Range("E25").Value = 12
$ws.Range("F25").Value = 10
$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 31
$ws.Range("E26").Value = 14
$ws.Range("F26").Value = 11
$ws.Range("B27").Value = 0
$ws.Range("C27").Value = 8
$ws.Range("D27").Value = 37
$ws.Range("E27").Value = 18
$ws.Range("F27").Value = 11
$ws.Range("B28").Value = 0
$ws.Range("C28").Value = 10
$ws.Range("D28").Value = 39
$ws.Range("E28").Value = 16
$ws.Range("F28").Value = 13
$ws.Range("B29").Value = 0
$ws.Range("C29").Value = 11
$ws.Range("D29").Value = 44
$ws.Range("E29").Value = 19
$ws.Range("F29").Value = 14
$ws.Range("B30").Value = 0
$ws.Range("C30").Value = 15
$ws.Range("D30").Value = 49
$ws.Range("E30").Value = 24
$ws.Range("F30").Value = 10
$ws.Range("B31").Value = 0
$ws.Range("C31").Value = 17
$ws.Range("D31").Value = 55
$ws.Range("E31").Value = 28
$ws.Range("F31").Value = 10
$ws.Range("B32").Value = 0
$ws.Range("C32").Value = 19
$ws.Range("D32").Value = 62
$ws.Range("E32").Value = 29
$ws.Range("F32").Value = 14
